$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 64959
$ws.Range("B2").Value = "Amanda Cavalcanti"
$ws.Range("C2").Value = "P&D"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45098
$ws.Range("G2").Value = 2619.46

# Row 3
$ws.Range("A3").Value = 12679
$ws.Range("B3").Value = "Thales Gonçalves"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45084
$ws.Range("G3").Value = 9022.33

# Row 4
$ws.Range("A4").Value = 77007
$ws.Range("B4").Value = "Thales Farias"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45081
$ws.Range("G4").Value = 7493.83

# Row 5
$ws.Range("A5").Value = 33219
$ws.Range("B5").Value = "Marina Freitas"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45087
$ws.Range("G5").Value = 10392.07

# Row 6
$ws.Range("A6").Value = 46484
$ws.Range("B6").Value = "Maria Vitória Gomes"
$ws.Range("C6").Value = "Jurídico"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45095
$ws.Range("G6").Value = 11813.77

# Row 7
$ws.Range("A7").Value = 65427
$ws.Range("B7").Value = "Dra. Ana Luiza Melo"
$ws.Range("C7").Value = "TI"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45105
$ws.Range("G7").Value = 6042.87

# Row 8
$ws.Range("A8").Value = 95347
$ws.Range("B8").Value = "Benício Vieira"
$ws.Range("C8").Value = "Operações"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45097
$ws.Range("G8").Value = 7483.32

# Row 9
$ws.Range("A9").Value = 14431
$ws.Range("B9").Value = "Alana Viana"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45089
$ws.Range("G9").Value = 11777.59

# Row 10
$ws.Range("A10").Value = 47822
$ws.Range("B10").Value = "Sr. Noah da Costa"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45094
$ws.Range("G10").Value = 4455.59

# Row 11
$ws.Range("A11").Value = 58429
$ws.Range("B11").Value = "Julia Nascimento"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45083
$ws.Range("G11").Value = 5524.46
